$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared string "Deferred to 0.0.11" shown (highlighted red) in several
#     rows' column C, using a new style (Arial font, solid red fill). ---
$deferredCells = @("C6", "C14", "C16", "C18", "C19", "C23", "C28")
foreach ($addr in $deferredCells) {
    $cell = $ws.Range($addr)
    $cell.Value = "Deferred to 0.0.11"
    $cell.Font.Name = "Arial"
    $cell.Interior.Color = 255
}

# --- Row 25: mark completed "Yes" plus a completion date/time ---
$ws.Range("C25").Value = "Yes"
$ws.Range("C4").Copy()
$ws.Range("C25").PasteSpecial(-4122)

$ws.Range("D25").Value = 45451.746527777781
$ws.Range("D3").Copy()
$ws.Range("D25").PasteSpecial(-4122)

# --- Row 26: mark completed "Yes" ---
$ws.Range("C26").Value = "Yes"
$ws.Range("C4").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Update the active cell selection to reflect where the edits were made ---
[void]$ws.Range("C14").Select()
